$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.420.58"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "1.566.42"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'209.28"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "'0.501"
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'22.02"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("E10").Value = "  +0.04%  "
$ws.Range("D11").Value = "'0.0867"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "1.789.92"
$ws.Range("E12").Value = "  -1.05%  "
$ws.Range("D13").Value = "1.563.49"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").Value = "'0.517"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").Value = "'63.53"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "27.396.32"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "'213.12"
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "'7.27"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "'9.55"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").Value = "'153.49"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'6.73"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'14.97"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").Value = "1.373.23"
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("D36").Value = "'0.965"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").Value = "'0.531"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'0.972"
$ws.Range("E42").Value = "  -0.79%  "
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("D44").Value = "'64.09"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("D47").Value = "1.701.90"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "'85.66"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("D49").Value = "0.0₇0988"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").Value = "'0.0955"
$ws.Range("E51").Value = "  -0.86%  "
